$d = $word.ActiveDocument

function New-XmlPackage($bodyXml) {
  return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr>'
$rPrB = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/></w:rPr>'

# ---------------------------------------------------------------------
# Edit 1: split the "On trip one / On trip two" sentence, wrapping the
# "one," and "two," words in grammar-check proofErr markers (as Word's
# grammar checker does when it flags a comma splice / fragment).
# ---------------------------------------------------------------------
$find = $d.Content
$find.Find.Execute("On trip one, take the parrot, leaving the cat and birdseed together. On trip two, bring the birdseed, leaving the cat alone. On the way back from trip two, bring the parrot back to the original side, leaving the birdseed. Before leaving for trip three, trade the parrot for the cat and leave the cat on the second side with the birdseed. Go back for the parrot, and make one more trip across. ") | Out-Null

$target = $d.Range($find.Start, $find.End)
$target.Text = ""
$insertPoint = $d.Range($find.Start, $find.Start)

$body = "<w:p>"
$body += "<w:r>$rPr<w:t xml:space=`"preserve`">On trip </w:t></w:r>"
$body += "<w:proofErr w:type=`"gramStart`"/>"
$body += "<w:r>$rPr<w:t>one,</w:t></w:r>"
$body += "<w:proofErr w:type=`"gramEnd`"/>"
$body += "<w:r>$rPr<w:t xml:space=`"preserve`"> take the parrot, leaving the cat and birdseed together. On trip </w:t></w:r>"
$body += "<w:proofErr w:type=`"gramStart`"/>"
$body += "<w:r>$rPr<w:t>two,</w:t></w:r>"
$body += "<w:proofErr w:type=`"gramEnd`"/>"
$body += "<w:r>$rPr<w:t xml:space=`"preserve`"> bring the birdseed, leaving the cat alone. On the way back from trip two, bring the parrot back to the original side, leaving the birdseed. Before leaving for trip three, trade the parrot for the cat and leave the cat on the second side with the birdseed. Go back for the parrot, and make one more trip across. </w:t></w:r>"
$body += "</w:p>"

$insertPoint.InsertXML((New-XmlPackage $body)) | Out-Null

# ---------------------------------------------------------------------
# Edit 2: split the "reaches 10, then 100, then 1,000" sentence,
# wrapping the second "then" in grammar-check proofErr markers.
# ---------------------------------------------------------------------
$find = $d.Content
$find.Find.Execute("The overall goal is to figure out which finger the girl will end on when she reaches 10, then 100, then 1,000. ") | Out-Null

$target = $d.Range($find.Start, $find.End)
$target.Text = ""
$insertPoint = $d.Range($find.Start, $find.Start)

$body = "<w:p>"
$body += "<w:r>$rPr<w:t xml:space=`"preserve`">The overall goal is to figure out which finger the girl will end on when she reaches 10, then 100, </w:t></w:r>"
$body += "<w:proofErr w:type=`"gramStart`"/>"
$body += "<w:r>$rPr<w:t>then</w:t></w:r>"
$body += "<w:proofErr w:type=`"gramEnd`"/>"
$body += "<w:r>$rPr<w:t xml:space=`"preserve`"> 1,000. </w:t></w:r>"
$body += "</w:p>"

$insertPoint.InsertXML((New-XmlPackage $body)) | Out-Null

# ---------------------------------------------------------------------
# Edit 3: fill in the "3) Identify potential solutions." section for
# the "Predicting Fingers" problem, which was previously left blank
# (just the lone _GoBack bookmark sitting in an empty paragraph).
# We replace that single empty paragraph with the heading paragraph
# (now bold, with its original paraId retained) plus six freshly
# authored paragraphs, moving the _GoBack bookmark to trail the very
# last bit of new text (mirroring where Word leaves it after a typing
# session).
# ---------------------------------------------------------------------
$find = $d.Content
$find.Find.Execute("3) Figure out what finger the girl will end on when she reaches 1,000.") | Out-Null
$blankPara = $find.Paragraphs(1).Next()
$target = $d.Range($blankPara.Range.Start, $blankPara.Range.End)

$body = ""
$body += '<w:p w14:paraId="48225126" w14:textId="77777777" w:rsidR="0033266B" w:rsidRPr="0033266B" w:rsidRDefault="0033266B" w:rsidP="0033266B">'
$body += "<w:pPr><w:spacing w:line=`"480`" w:lineRule=`"auto`"/><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/><w:b/></w:rPr></w:pPr>"
$body += "<w:r>$rPrB<w:t>3) Identify potential solutions.</w:t></w:r>"
$body += "</w:p>"

$body += "<w:p>"
$body += "<w:pPr><w:spacing w:line=`"480`" w:lineRule=`"auto`"/><w:ind w:left=`"720`"/><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/><w:b/></w:rPr></w:pPr>"
$body += "<w:r>$rPrB<w:t>a) For each of the sub-problems you’ve discussed in #2, what is a possible solution?</w:t></w:r>"
$body += "</w:p>"

$body += "<w:p>"
$body += "<w:pPr><w:spacing w:line=`"480`" w:lineRule=`"auto`"/><w:ind w:left=`"720`"/><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/></w:rPr></w:pPr>"
$body += "<w:r>$rPrB<w:tab/></w:r>"
$body += "<w:r>$rPr<w:t>1) Count on your own hand. She will end on her first finger.</w:t></w:r>"
$body += "</w:p>"

$body += "<w:p>"
$body += "<w:pPr><w:spacing w:line=`"480`" w:lineRule=`"auto`"/><w:ind w:left=`"720`"/><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/></w:rPr></w:pPr>"
$body += "<w:r>$rPr<w:tab/><w:t>2) Figure out a pattern. The 10s will always land on the first or ring finger.</w:t></w:r>"
$body += "</w:p>"

$body += "<w:p>"
$body += "<w:pPr><w:spacing w:line=`"480`" w:lineRule=`"auto`"/><w:ind w:left=`"2160`"/><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/></w:rPr></w:pPr>"
$body += "<w:r>$rPr<w:t xml:space=`"preserve`">10 </w:t></w:r>"
$body += "<w:proofErr w:type=`"gramStart`"/>"
$body += "<w:r>$rPr<w:t>is</w:t></w:r>"
$body += "<w:proofErr w:type=`"gramEnd`"/>"
$body += "<w:r>$rPr<w:t xml:space=`"preserve`"> on the first finger, 20 and 30 are on the ring finger, then it continues by twos. 100 will land on her ring finger.</w:t></w:r>"
$body += "</w:p>"

$body += "<w:p>"
$body += "<w:pPr><w:spacing w:line=`"480`" w:lineRule=`"auto`"/><w:ind w:left=`"1440`"/><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/></w:rPr></w:pPr>"
$body += "<w:r>$rPr<w:t>3) Continue the pattern. When counting by 100s, switch off between the</w:t></w:r>"
$body += "</w:p>"

$body += "<w:p>"
$body += "<w:pPr><w:spacing w:line=`"480`" w:lineRule=`"auto`"/><w:ind w:left=`"1440`" w:firstLine=`"720`"/><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/></w:rPr></w:pPr>"
$body += "<w:proofErr w:type=`"gramStart`"/>"
$body += "<w:r>$rPr<w:t>first</w:t></w:r>"
$body += "<w:proofErr w:type=`"gramEnd`"/>"
$body += "<w:r>$rPr<w:t xml:space=`"preserve`"> and ring fingers. 1,000 will end on her first finger.</w:t></w:r>"
$body += '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$body += "</w:p>"

$target.InsertXML((New-XmlPackage $body)) | Out-Null

Write-Output "All edits applied"
